# Update "合肥-漫展信息" workbook:
#  - Remove the "合肥·乐帮•崩原铁绝only同人首展" event (2024-10-04) from the
#    "展览" (Exhibitions) sheet and the "全部类型" (All types) sheet — it no
#    longer shows up in the refreshed data pull.
#  - Refresh a handful of "想去人数" (interest-count) figures that ticked up
#    between pulls, and flip one event's "最低票价" (lowest price) between
#    "不可售" (not for sale) and a concrete number as ticket sales opened.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "展览" (Exhibitions) sheet
# ---------------------------------------------------------------------
$wsEx = $wb.Worksheets.Item("展览")

# Row 3 holds "合肥·乐帮•崩原铁绝only同人首展" (2024-10-04) — delete it;
# everything below shifts up one row.
$wsEx.Rows.Item(3).EntireRow.Delete()

# Renumber the leading index column (A) so it stays 0..N sequential.
$usedEx = $wsEx.UsedRange
$lastRowEx = $usedEx.Row + $usedEx.Rows.Count - 1
for ($r = 2; $r -le $lastRowEx; $r++) {
    $wsEx.Cells.Item($r, 1).Value = $r - 1
}

# Post-shift value refresh.
$wsEx.Range("F2").Value = 603       # 合肥·Holic动漫游戏展             601 -> 603
$wsEx.Range("G4").Value = "不可售"  # 合肥·首届火影忍者同人only        29.9 -> 不可售
$wsEx.Range("F5").Value = 361       # 合肥·W·A第五人格同人only2.0      360 -> 361
$wsEx.Range("G5").Value = 68        # 合肥·W·A第五人格同人only2.0      不可售 -> 68
$wsEx.Range("F6").Value = 1867      # 合肥·第九届环形宇宙动漫游戏嘉年华 1835 -> 1867
$wsEx.Range("F7").Value = 101       # 合肥·MAX特摄同人only2.0          100 -> 101

# ---------------------------------------------------------------------
# "全部类型" (All types) sheet
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Rows.Item(3).EntireRow.Delete()

$usedAll = $wsAll.UsedRange
$lastRowAll = $usedAll.Row + $usedAll.Rows.Count - 1
for ($r = 2; $r -le $lastRowAll; $r++) {
    $wsAll.Cells.Item($r, 1).Value = $r - 1
}

$wsAll.Range("F2").Value = 603       # 合肥·Holic动漫游戏展             601 -> 603
$wsAll.Range("G4").Value = "不可售"  # 合肥·首届火影忍者同人only        29.9 -> 不可售
$wsAll.Range("F5").Value = 361       # 合肥·W·A第五人格同人only2.0      360 -> 361
$wsAll.Range("G5").Value = 68        # 合肥·W·A第五人格同人only2.0      不可售 -> 68
$wsAll.Range("F10").Value = 1867     # 合肥·第九届环形宇宙动漫游戏嘉年华 1835 -> 1867
$wsAll.Range("F11").Value = 101      # 合肥·MAX特摄同人only2.0          100 -> 101
